$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Npc")
$ws.Range("A1").Value = "test"
